# Chapitre 3 - Structures conditionnelles
# Feat(chp3): Ajout des commentaire de documentation
#
# Slide 4 (sldId 287) : on enrichit les deux explications sur les
# commentaires Python (commentaire sur une ligne / commentaire multilignes)
# avec une precision entre parentheses, et on repositionne/redimensionne
# les zones de texte et les images de capture d'ecran en consequence.
#
# NB: PowerPoint's Shape.Left/Top/Width/Height are Single-precision
# (32-bit float) point values under the hood. To land on the exact target
# EMU after the point -> EMU round-trip, the literals below are the
# closest representable point values (pre-computed) rather than the naive
# EMU/12700 division.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# --- ZoneTexte 8 (id 9) : "Les commentaires sur une ligne: ..." ---
$shComment1 = Get-ShapeById $s 9
$shComment1.Left = 27.627717971801758      # 350872 EMU
$shComment1.Top = 142.6991424560547        # 1812279 EMU
$shComment1.Width = 954.0694580078125      # 12116682 EMU
$shComment1.Height = 55.73905563354492     # 707886 EMU

$para1 = $shComment1.TextFrame.TextRange.Paragraphs(1)
$para1.Runs(2).Text = "commence par le caractère dièse # (utilisé pour commenter le code utile pour le développeur.)."

# --- ZoneTexte 9 (id 10) : "Les commentaires multilignes: ..." ---
$shComment2 = Get-ShapeById $s 10
$shComment2.Left = 39.44346618652344        # 500932 EMU
$shComment2.Top = 290.9941101074219         # 3695625 EMU

$para2 = $shComment2.TextFrame.TextRange.Paragraphs(1)
$para2.Runs(2).Text = "commentaires sur plusieurs lignes encadrées par 3 guillemets doubles ou simples (utile pour faire de la documentation de code pour les autres.)."

# --- Image 12 (id 13) : capture d'ecran du commentaire sur une ligne ---
$shImg1 = Get-ShapeById $s 13
$shImg1.Top = 205.15867614746094            # 2605515 EMU

# --- Image 14 (id 15) : capture d'ecran du commentaire multilignes ---
$shImg2 = Get-ShapeById $s 15
$shImg2.Left = 44.660552978515625           # 567189 EMU
$shImg2.Top = 388.9093017578125             # 4939148 EMU

# --- Image 16 (id 17) : seconde capture d'ecran du commentaire multilignes ---
$shImg3 = Get-ShapeById $s 17
$shImg3.Left = 497.5592346191406            # 6319002 EMU
$shImg3.Top = 388.9093017578125             # 4939148 EMU
